$d = $word.ActiveDocument

# --- Table 3: update mean friction factor / Reynolds number values ---
$d.Content.Find.Execute("0.0439982", $true, $false, $false, $false, $false, $true, 1, $false, "0.0440", 2)
$d.Content.Find.Execute("11187.73", $true, $false, $false, $false, $false, $true, 1, $false, "11200", 2)
$d.Content.Find.Execute("0.0383047", $true, $false, $false, $false, $false, $true, 1, $false, "0.0383", 2)
$d.Content.Find.Execute("18646.21", $true, $false, $false, $false, $false, $true, 1, $false, "18600", 2)
$d.Content.Find.Execute("0.0347399", $true, $false, $false, $false, $false, $true, 1, $false, "0.0347", 2)
$d.Content.Find.Execute("26104.70", $true, $false, $false, $false, $false, $true, 1, $false, "26100", 2)

# --- Table 4: update percent difference value ---
$d.Content.Find.Execute("24.1%", $true, $false, $false, $false, $false, $true, 1, $false, "23.9%", 2)

# --- Insert a Reynolds-number question (bracketed by horizontal rules) just
#     before the "We see that the lower flow rates..." paragraph, and switch
#     that paragraph's style from BodyText to FirstParagraph. ---
$hr = '<w:p><w:r><w:pict><v:rect style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t"/></w:pict></w:r></w:p>'

$r = $d.Content
$r.Find.Execute("We see that the lower flow rates yield lower Reynolds numbers and higher friction factors, consistent with expectations from the Moody chart.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $r.Paragraphs(1).Range

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office">
<w:body>
$hr
<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">My reynolds numbers are off compared to the original report, are my reynolds numbers correct?</w:t></w:r></w:p>
$hr
<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">We see that the lower flow rates yield lower Reynolds numbers and higher friction factors, consistent with expectations from the Moody chart.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
$target.InsertXML($xml)
